# Edit script: rename "Non-Targeted Analysis" / "Non-targeted analysis" (NTA)
# terminology to "Suspect Screening Analysis" / "suspect screening analysis"
# across the poster title, Methods blurb and Figure 1 caption, and refresh the
# auto date placeholders (datetimeFigureOut fields) on the slide master and
# every slide layout from 8/22/2023 to 4/4/2024.

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

function Set-DatePlaceholderText($shapes, [string]$newText) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "*Date*") {
            $tr = $sh.TextFrame.TextRange
            $tr.Text = $newText
        }
    }
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. Poster title (shape id 8, "TextBox 7"):
#    "...Using Non-Targeted Analysis" -> "...Using Suspect Screening Analysis"
#    Font shrinks 82pt -> 78pt and the autosized text box shrinks accordingly.
# ---------------------------------------------------------------------
$shTitle = Get-ShapeById $s.Shapes 8
$trTitle = $shTitle.TextFrame.TextRange
$trTitle.Text = "Developing Chemical Signatures for 5 Categories of Household Products Using Suspect Screening Analysis"
$trTitle.Font.Size = 78
$shTitle.Height = 196.2984251968504

# ---------------------------------------------------------------------
# 2. Methods bullet (shape id 22, "TextBox 21"):
#    "Non-targeted analysis (NTA) using ..." -> "Suspect screening analysis (NTA) using ..."
# ---------------------------------------------------------------------
$shMethods = Get-ShapeById $s.Shapes 22
$trMethods = $shMethods.TextFrame.TextRange
$fullMethods = $trMethods.Text
$oldMethods = "Non-targeted"
$idxMethods = $fullMethods.IndexOf($oldMethods)
$subMethods = $trMethods.Characters($idxMethods + 1, $oldMethods.Length)
$subMethods.Text = "Suspect screening"

# ---------------------------------------------------------------------
# 3. Figure 1 caption (shape id 40, "TextBox 39"):
#    "Workflow of non-targeted analysis ..." -> "Workflow of suspect screening analysis ..."
# ---------------------------------------------------------------------
$shFig1 = Get-ShapeById $s.Shapes 40
$trFig1 = $shFig1.TextFrame.TextRange
$fullFig1 = $trFig1.Text
$oldFig1 = "Workflow of non-targeted analysis of products from 5 types of household consumer products. Products were extracted with dichloromethane (DCM). After addition of an internal standard, each extraction was analyzed via GC X GC-TOFMS to obtain its mass spectra. The spectra were matched to the 2017 NIST database and analytical standards were used to confirm a subset of the chemical identifications. Chemicals were annotated by reported or predicted functional uses"
$newFig1 = "Workflow of suspect screening analysis of products from 5 types of household consumer products. Products were extracted with dichloromethane (DCM). After addition of an internal standard, each extraction was analyzed via GC X GC-TOFMS to obtain its mass spectra. The spectra were matched to the 2017 NIST database and analytical standards were used to confirm a subset of the chemical identifications. Chemicals were annotated by reported or predicted functional uses"
$idxFig1 = $fullFig1.IndexOf($oldFig1)
$subFig1 = $trFig1.Characters($idxFig1 + 1, $oldFig1.Length)
$subFig1.Text = $newFig1

# ---------------------------------------------------------------------
# 4. Refresh the "datetimeFigureOut" date placeholders (slide master +
#    every custom layout) from 8/22/2023 to 4/4/2024.
# ---------------------------------------------------------------------
$design = $p.Designs.Item(1)
$slideMaster = $design.SlideMaster

Set-DatePlaceholderText $slideMaster.Shapes "4/4/2024"

$layouts = $slideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Set-DatePlaceholderText $layouts.Item($i).Shapes "4/4/2024"
}
